$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New rows of data for studia "S1-STA": pesel, wynik, zakwalifikowany, przyjety
$data = @(
    @(98090100001, 100, 1, 1),
    @(98090100002, 99, 1, 1),
    @(98090100003, 98, 1, 1),
    @(98090100004, 97, 1, 1),
    @(98090100005, 96, 1, 1),
    @(98090100006, 95, 1, 1),
    @(98090100007, 94, 1, 1),
    @(98090100008, 93, 1, 1),
    @(98090100009, 92, 1, 1),
    @(98090100010, 91, 1, 1),
    @(98090100011, 90, 1, 0),
    @(98090100012, 89, 1, 0),
    @(98090100013, 88, 1, 0),
    @(98090100014, 87, 1, 0),
    @(98090100015, 86, 1, 0),
    @(98090100016, 85, 1, 0),
    @(98090100017, 84, 1, 0),
    @(98090100018, 83, 1, 0),
    @(98090100019, 82, 1, 0),
    @(98090100020, 81, 0, 0),
    @(98090100021, 80, 0, 0),
    @(98090100022, 79, 0, 0),
    @(98090100023, 78, 0, 0),
    @(98090100024, 77, 0, 0),
    @(98090100025, 76, 0, 0),
    @(98090100026, 75, 0, 0),
    @(98090100027, 74, 0, 0),
    @(98090100028, 73, 0, 0),
    @(98090100029, 72, 0, 0),
    @(98090100030, 71, 0, 0)
)

$startRow = 145
for ($i = 0; $i -lt $data.Length; $i++) {
    $row = $startRow + $i
    $rowData = $data[$i]
    $ws.Cells.Item($row, 1).Value = $rowData[0]
    $ws.Cells.Item($row, 2).Value = "S1-STA"
    $ws.Cells.Item($row, 4).Value = $rowData[1]
    $ws.Cells.Item($row, 7).Value = $rowData[2]
    $ws.Cells.Item($row, 8).Value = $rowData[3]
}

# Update the selection/view to match the new data range
$ws.Range("D145:D174").Select()
